$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 2
    4  = 0
    5  = 2
    6  = 0
    7  = 7
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 3
    14 = 0
    15 = 0
    16 = 0
    17 = 2
    19 = 1
    20 = 2
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 2
    27 = 1
    28 = 2
    29 = 3
    30 = 1
    31 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
